$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.1507537688442211
$ws.Cells.Item(2, 3).Value = 0.6407035175879398
$ws.Cells.Item(2, 10).Value = 0.01256281407035176
$ws.Cells.Item(2, 16).Value = 0.1331658291457286
$ws.Cells.Item(2, 19).Value = 0.06281407035175879
$ws.Cells.Item(3, 2).Value = 0.007067137809187279
$ws.Cells.Item(3, 3).Value = 0.02473498233215548
$ws.Cells.Item(3, 10).Value = 0.03180212014134275
$ws.Cells.Item(3, 16).Value = 0.7667844522968198
$ws.Cells.Item(3, 19).Value = 0.1696113074204947
$ws.Cells.Item(4, 10).Value = 0.01204819277108434
$ws.Cells.Item(4, 16).Value = 0.7108433734939759
$ws.Cells.Item(4, 19).Value = 0.2771084337349398
$ws.Cells.Item(6, 2).Value = 0.07964601769911504
$ws.Cells.Item(6, 4).Value = 0.02359882005899705
$ws.Cells.Item(6, 6).Value = 0.05309734513274336
$ws.Cells.Item(6, 10).Value = 0.2418879056047198
$ws.Cells.Item(6, 15).Value = 0.02064896755162242
$ws.Cells.Item(6, 17).Value = 0.1799410029498525
$ws.Cells.Item(6, 18).Value = 0.08259587020648967
$ws.Cells.Item(6, 19).Value = 0.3185840707964602
$ws.Cells.Item(7, 2).Value = 0.1056338028169014
$ws.Cells.Item(7, 4).Value = 0.04577464788732395
$ws.Cells.Item(7, 6).Value = 0.08450704225352113
$ws.Cells.Item(7, 10).Value = 0.09154929577464789
$ws.Cells.Item(7, 15).Value = 0.0176056338028169
$ws.Cells.Item(7, 17).Value = 0.2253521126760563
$ws.Cells.Item(7, 18).Value = 0.09507042253521127
$ws.Cells.Item(7, 19).Value = 0.3345070422535211
$ws.Cells.Item(8, 2).Value = 0.09781021897810219
$ws.Cells.Item(8, 4).Value = 0.01605839416058394
$ws.Cells.Item(8, 6).Value = 0.07153284671532846
$ws.Cells.Item(8, 10).Value = 0.1109489051094891
$ws.Cells.Item(8, 15).Value = 0.0218978102189781
$ws.Cells.Item(8, 17).Value = 0.2102189781021898
$ws.Cells.Item(8, 18).Value = 0.08759124087591241
$ws.Cells.Item(8, 19).Value = 0.3839416058394161
$ws.Cells.Item(9, 2).Value = 0.1446280991735537
$ws.Cells.Item(9, 4).Value = 0.02892561983471074
$ws.Cells.Item(9, 5).Value = 0.004132231404958678
$ws.Cells.Item(9, 6).Value = 0.07024793388429752
$ws.Cells.Item(9, 10).Value = 0.07851239669421488
$ws.Cells.Item(9, 15).Value = 0.01239669421487603
$ws.Cells.Item(9, 17).Value = 0.2355371900826446
$ws.Cells.Item(9, 18).Value = 0.05371900826446281
$ws.Cells.Item(9, 19).Value = 0.371900826446281
$ws.Cells.Item(10, 2).Value = 0.08980454305335446
$ws.Cells.Item(10, 4).Value = 0.02377179080824089
$ws.Cells.Item(10, 6).Value = 0.07237189646064449
$ws.Cells.Item(10, 10).Value = 0.109878499735869
$ws.Cells.Item(10, 15).Value = 0.02218700475435816
$ws.Cells.Item(10, 17).Value = 0.2282091917591125
$ws.Cells.Item(10, 18).Value = 0.09561542525092445
$ws.Cells.Item(10, 19).Value = 0.358161648177496
$ws.Cells.Item(11, 7).Value = 0.1334894613583138
$ws.Cells.Item(11, 10).Value = 0.107728337236534
$ws.Cells.Item(11, 11).Value = 0.1826697892271663
$ws.Cells.Item(11, 12).Value = 0.5667447306791569
$ws.Cells.Item(11, 19).Value = 0.00936768149882904
$ws.Cells.Item(12, 7).Value = 0.7579365079365079
$ws.Cells.Item(12, 10).Value = 0.1706349206349206
$ws.Cells.Item(12, 11).Value = 0.007936507936507936
$ws.Cells.Item(12, 12).Value = 0.03571428571428571
$ws.Cells.Item(12, 19).Value = 0.02777777777777778
$ws.Cells.Item(13, 7).Value = 0.7254901960784313
$ws.Cells.Item(13, 10).Value = 0.196078431372549
$ws.Cells.Item(13, 19).Value = 0.07843137254901961
$ws.Cells.Item(14, 7).Value = 0.8
$ws.Cells.Item(14, 10).Value = 0.2
$ws.Cells.Item(15, 6).Value = 0.01449275362318841
$ws.Cells.Item(15, 8).Value = 0.1884057971014493
$ws.Cells.Item(15, 9).Value = 0.05797101449275362
$ws.Cells.Item(15, 10).Value = 0.3304347826086956
$ws.Cells.Item(15, 11).Value = 0.06666666666666667
$ws.Cells.Item(15, 13).Value = 0.005797101449275362
$ws.Cells.Item(15, 15).Value = 0.04057971014492753
$ws.Cells.Item(15, 19).Value = 0.2956521739130435
$ws.Cells.Item(16, 6).Value = 0.02531645569620253
$ws.Cells.Item(16, 8).Value = 0.1772151898734177
$ws.Cells.Item(16, 9).Value = 0.06012658227848101
$ws.Cells.Item(16, 10).Value = 0.3955696202531646
$ws.Cells.Item(16, 11).Value = 0.120253164556962
$ws.Cells.Item(16, 13).Value = 0.02848101265822785
$ws.Cells.Item(16, 15).Value = 0.05696202531645569
$ws.Cells.Item(16, 19).Value = 0.1360759493670886
$ws.Cells.Item(17, 6).Value = 0.01703800786369594
$ws.Cells.Item(17, 8).Value = 0.1939711664482307
$ws.Cells.Item(17, 9).Value = 0.07077326343381389
$ws.Cells.Item(17, 10).Value = 0.4338138925294889
$ws.Cells.Item(17, 11).Value = 0.08781127129750983
$ws.Cells.Item(17, 13).Value = 0.01310615989515072
$ws.Cells.Item(17, 14).Value = 0.002621231979030144
$ws.Cells.Item(17, 15).Value = 0.07863695937090433
$ws.Cells.Item(17, 19).Value = 0.1022280471821756
$ws.Cells.Item(18, 6).Value = 0.02564102564102564
$ws.Cells.Item(18, 8).Value = 0.1826923076923077
$ws.Cells.Item(18, 9).Value = 0.0608974358974359
$ws.Cells.Item(18, 10).Value = 0.4551282051282051
$ws.Cells.Item(18, 11).Value = 0.09935897435897435
$ws.Cells.Item(18, 13).Value = 0.009615384615384616
$ws.Cells.Item(18, 15).Value = 0.04166666666666666
$ws.Cells.Item(18, 19).Value = 0.125
$ws.Cells.Item(19, 6).Value = 0.0137299771167048
$ws.Cells.Item(19, 8).Value = 0.2053775743707094
$ws.Cells.Item(19, 9).Value = 0.07608695652173914
$ws.Cells.Item(19, 10).Value = 0.3958810068649886
$ws.Cells.Item(19, 11).Value = 0.1052631578947368
$ws.Cells.Item(19, 13).Value = 0.01659038901601831
$ws.Cells.Item(19, 14).Value = 0.002288329519450801
$ws.Cells.Item(19, 15).Value = 0.07608695652173914
$ws.Cells.Item(19, 19).Value = 0.108695652173913
